# BatxDetail.xlsx - "Multiple URS Updates; GenTables updates L1, L3;
# Replies to SKL inquiries"
#
# The ProcNote sheet (GenTable field dictionary) gets two new rows
# documenting two new DB columns used elsewhere in the L4 batch job
# table layouts:
#   9  RepayBank   / 扣款銀行         / VARCHAR2  / 3
#   10 PayIntDate  / 銀扣期款應繳日    / Decimald  / 8
#
# The author was last looking at the ProcNote sheet (tab + selection),
# so we finish by activating that sheet and selecting H26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProcNote")

# --- New row 14 (item 9: RepayBank) -----------------------------------
# Clone formatting from the row above (row 12), which already has the
# plain (non-shaded) style used for odd-numbered field rows.
$ws.Range("B12:H12").Copy()
$ws.Range("B14:H14").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B14").Value = 9
$ws.Range("C14").Value = "RepayBank"
$ws.Range("D14").Value = "扣款銀行"
$ws.Range("E14").Value = "VARCHAR2"
$ws.Range("F14").Value = 3

# --- New row 15 (item 10: PayIntDate) ---------------------------------
# Clone formatting from row 13, which carries the special H-column
# style used for the table's last data row.
$ws.Range("B13:H13").Copy()
$ws.Range("B15:H15").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B15").Value = 10
$ws.Range("C15").Value = "PayIntDate"
$ws.Range("D15").Value = "銀扣期款應繳日"
$ws.Range("E15").Value = "Decimald"
$ws.Range("F15").Value = 8

$excel.CutCopyMode = $false

# --- Active sheet / selection bookkeeping -----------------------------
$ws.Activate() | Out-Null
$ws.Range("H26").Select() | Out-Null
